# Helper: force a run boundary at an absolute document character position
# by briefly wrapping it with a throw-away bookmark (Word splits runs at
# bookmark boundaries; the split survives removal of the bookmark itself).
function Split-RunAt($doc, $pos, $tag) {
    $r = $doc.Range($pos, $pos)
    $bmName = "TmpSplit" + $tag
    $doc.Bookmarks.Add($bmName, $r) | Out-Null
    $doc.Bookmarks.Item($bmName).Delete()
}

$d = $word.ActiveDocument

# =====================================================================
# Hunk 1: paragraph "6 - Publicar repositório remoto PUSH." -- collapse
# the five separate runs ("6 - " / "Publicar repositório " / "remoto " /
# "PUSH" / ".") into a single run, leaving the following run
# (" (Após um ...") untouched.
# =====================================================================
$p6 = $d.Paragraphs.Item(8)
$full = $p6.Range
$full.Find.Execute("6 - Publicar repositório remoto PUSH.") | Out-Null
$full.Delete()
$full.InsertBefore("6 - Publicar repositório remoto PUSH.")

# =====================================================================
# Hunk 2: after the (former) bookmark-only paragraph, insert three new
# paragraphs, then one more empty paragraph, and relocate the _GoBack
# bookmark into the last of the new paragraphs.
# =====================================================================
$pA = $d.Paragraphs.Item(9)

# Create three new empty paragraphs right after $pA (they will become
# paragraphs B, C, D), plus a trailing empty paragraph (E).
$pA.Range.InsertParagraphAfter()
$pA.Range.InsertParagraphAfter()
$pA.Range.InsertParagraphAfter()
$pA.Range.InsertParagraphAfter()

$pB = $d.Paragraphs.Item(10)
$pC = $d.Paragraphs.Item(11)
$pD = $d.Paragraphs.Item(12)

# ---- Paragraph B: "Continuando o trabalho armazenado em um repositório
# remoto" + " (clone)" (two runs) ----
$pB.Range.InsertBefore("Continuando o trabalho armazenado em um repositório remoto (clone)")
$pB2 = $d.Paragraphs.Item(10)
$splitPos = $pB2.Range.Start + 60
Split-RunAt $d $splitPos "B"
Write-Output "pB done: [$($d.Paragraphs.Item(10).Range.Text)]"

# ---- Paragraph C: "1 – Abrir GitHub Desktop" (single run) ----
$pC2 = $d.Paragraphs.Item(11)
$pC2.Range.InsertBefore("1 – Abrir GitHub Desktop")
Write-Output "pC done: [$($d.Paragraphs.Item(11).Range.Text)]"

# ---- Paragraph D: "2 " | "–" | " " | "Em arquivo, selecionar clone
# repository" + bookmark(_GoBack) + ".." ----
$pD2 = $d.Paragraphs.Item(12)
$pD2.Range.InsertBefore("2 – Em arquivo, selecionar clone repositoryX..")
$pD3 = $d.Paragraphs.Item(12)
$dStart = $pD3.Range.Start
$split1 = $dStart + 2
$split2 = $dStart + 3
$split3 = $dStart + 4
Split-RunAt $d $split1 "D1"
Split-RunAt $d $split2 "D2"
Split-RunAt $d $split3 "D3"
Write-Output "pD done: [$($d.Paragraphs.Item(12).Range.Text)]"

# Relocate the _GoBack bookmark to sit right after "...repository" and
# before "..": wrap the placeholder "X" char with _GoBack, then delete it.
$pD4 = $d.Paragraphs.Item(12)
$xPos = $dStart + 44
$xRange = $d.Range($xPos, $xPos + 1)
Write-Output "xRange text: [$($xRange.Text)]"
$d.Bookmarks.Add("_GoBack", $xRange) | Out-Null
$xRange2 = $d.Range($xPos, $xPos + 1)
$xRange2.Text = ""
Write-Output "pD final: [$($d.Paragraphs.Item(12).Range.Text)]"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($p.Range.Text)]"
}
